$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)

        $v1 = $c1.Value2
        $v2 = $c2.Value2

        # Skip cells that already hold the same value on both rows so we
        # don't needlessly rewrite (and reformat) untouched data.
        if ($v1 -ceq $v2) {
            continue
        }

        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Column A holds the fixed sequential row id and must stay untouched.
# Columns B..AD (2..30) carry the actual match data that needs to swap
# between the two rows in each pair.
Swap-RowData 168 169 2 30
Swap-RowData 180 181 2 30
Swap-RowData 184 185 2 30
